# "atualizei dados bibi e add" - refresh the vendas_atipicas sample:
#  - Drop the four now-obsolete rows (2025-06-12 x2, 2025-06-13, 2025-06-14)
#    that used to sit at the top of the table.
#  - That leaves the former rows 6-10 (2025-06-16 .. 2025-06-24) shifted up
#    into rows 2-6, already carrying the right date/client/id_venda/produto
#    text - only their estoque_atualizado (G) and desvio_padrao (I) figures
#    needed refreshing to the latest numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the four obsolete data rows; rows 6-10 shift up to become rows 2-6.
$ws.Range("A2:I5").EntireRow.Delete() | Out-Null

# Refresh the numeric columns that changed for the surviving rows.
$ws.Cells.Item(2, 7).Value = -440
$ws.Cells.Item(2, 9).Value = 0.28

$ws.Cells.Item(3, 7).Value = -295

$ws.Cells.Item(4, 7).Value = -440
$ws.Cells.Item(4, 9).Value = 0.28

$ws.Cells.Item(5, 7).Value = -440
$ws.Cells.Item(5, 9).Value = 0.28

# Row 6 (2025-06-24 / SMART WATCH HMASTON INK12) is unchanged.
